$d = $word.ActiveDocument

# --- Edit 1: MCQ answer text split -------------------------------------
# "...your page design.")?"  ->  "...your page design method.")?"
# (a new word " method" is inserted right before the closing '."')
$rdquo = [char]0x201D
$r = $d.Content
$found = $r.Find.Execute(
    "would Actually use these methods depends on your page design." + $rdquo,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "would Actually use these methods depends on your page design method." + $rdquo,
    2)

# --- Edit 2: "case 1:" -> "case 0:" -------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    "    case 1:", $true, $false, $false, $false, $false, $true, 1, $false,
    "    case 0:", 2)

# --- Edit 3: "case 2:" -> "case 1:" -------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute(
    "    case 2:", $true, $false, $false, $false, $false, $true, 1, $false,
    "    case 1:", 2)

Write-Output "edit1=$found edit2=$found2 edit3=$found3"
